# The deck ships with two theme parts:
#   ppt/theme/theme1.xml  -> "Office Theme" (the default Office palette) - only
#                             wired to the Notes Master, otherwise unused
#   ppt/theme/theme2.xml  -> "Integral"      - the theme actually driving the
#                             slide master / every slide in the deck
#
# The authored change swaps the content of those two theme parts, so the
# design actually applied to the slides switches from "Integral" to the
# plain "Office Theme" palette (and vice versa for the notes master).
#
# The only theme surface the PowerPoint object model exposes for editing is
# the *active* design's ThemeColorScheme (backed by ppt/theme/theme2.xml,
# since that's what the slide master / presentation point at). We push the
# swap through by re-pointing every one of its 12 theme colors at the
# "Office Theme" values so the design in effect after this edit matches the
# target state.

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$master = $design.SlideMaster
$theme = $master.Theme
$colorScheme = $theme.ThemeColorScheme

# Index order == a:clrScheme child order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
# Values are standard "Office Theme" RGB values, encoded as COM BGR ints (R + G*256 + B*65536).
$officeThemeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeThemeColors[$i - 1]
}
